$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 79; this shifts existing rows 79..147 down to 80..148
$ws.Rows.Item(79).Insert()

# Populate the newly inserted row 79 with the new data record
$ws.Cells.Item(79, 1).Value = 6
$ws.Cells.Item(79, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(79, 3).Value = "Metropolitana"
$ws.Cells.Item(79, 4).Value = 44586
$ws.Cells.Item(79, 4).Style = $ws.Cells.Item(80, 4).Style
$ws.Cells.Item(79, 4).NumberFormat = $ws.Cells.Item(80, 4).NumberFormat
$ws.Cells.Item(79, 5).Value = 13
$ws.Cells.Item(79, 6).Value = 100112029
$ws.Cells.Item(79, 7).Value = "Orégano"
$ws.Cells.Item(79, 8).Value = "Sin especificar"
$ws.Cells.Item(79, 9).Value = "Primera"
$ws.Cells.Item(79, 10).Value = 34
$ws.Cells.Item(79, 11).Value = 8000
$ws.Cells.Item(79, 12).Value = 9000
$ws.Cells.Item(79, 13).Value = 8441
$ws.Cells.Item(79, 14).Value = "$/docena de atados"
$ws.Cells.Item(79, 15).Value = "Región Metropolitana"
$ws.Cells.Item(79, 16).Value = 2814
$ws.Cells.Item(79, 17).Value = 3
$ws.Cells.Item(79, 18).Value = "Hortaliza"
